$block3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I used a smart pointer to manage the </w:t></w:r><w:r><w:t xml:space="preserve">main </w:t></w:r><w:r><w:t>graph object in main.cpp.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">All graph files can be opened with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Nodepad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$block15 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I’m going to start by building the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GraphADT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class from the textbook (pg. </w:t></w:r><w:r><w:t>387)</w:t></w:r></w:p>
'@
$block16 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Now that the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GraphADT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class</w:t></w:r><w:r><w:t xml:space="preserve"> has been built, </w:t></w:r><w:r><w:t>I’m going to use the author’s adjacency matrix implementation (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>389)</w:t></w:r><w:r><w:t xml:space="preserve">. I did make a few changes for naming conventions. Namely, changing </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>numEdge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>numVertext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>numEdges</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>numVertices</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Also, I use v1 and v2 where the author uses </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and j.</w:t></w:r></w:p>
'@
$block29 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Write main that tests reading graph from file, creating three </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from three different vertices, and prints the graph to a file (as well as console)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"></w:p>
'@
$blockWritefn = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Write function for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$block42 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I’m </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gonna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> use a smart pointer </w:t></w:r><w:r><w:t>to store the graph in main.cpp</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Now that everything is set up, I’m </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gonna</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> write a function to generate a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from a given vertex that will return a pointer to the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> using Prim’s Algorithm f</w:t></w:r><w:r><w:t>rom the textbook.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">In my implementation of Prim’s Algorithm, </w:t></w:r><w:r><w:t xml:space="preserve">I’ve used the author’s source code and helper functions as a guide. </w:t></w:r><w:r><w:t xml:space="preserve">I had to look at them a lot before I understood what was going on. The way he wrote out the algorithm seemed much simpler than the code I was reading, and I still can’t match every bit one to one but I have enough understanding of what it is doing to modify it in order to accomplish my goals. I did have to edit the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>addEdgetoMST</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> part because for whatever reason, in my code it was only marking the edges directionally. </w:t></w:r><w:r><w:t xml:space="preserve">I accounted for this by simply mirroring the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>setEdge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> code. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I additionally had a problem where if I did any of the vertices’ MST by itself, it worked fine, but if I tried to add more than one to main it would crash. </w:t></w:r><w:r><w:t xml:space="preserve">I didn’t realize at first, but this was because I was never marking all the nodes as unvisited again after one of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> attempts ran.</w:t></w:r></w:p>
'@
$block44 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>I’ve built the read, write, and print functions and all of them work as expected.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I built the MCST function, which took a lot longer than expected, and got it to work in accordance with the lab directions. </w:t></w:r><w:r><w:t xml:space="preserve">It seems as though all three starting vertices give the same </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mcst</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, and I’m not sure if that is because of the graph that I’m using or a flaw in the logic I implemented.</w:t></w:r><w:r><w:t xml:space="preserve"> Because my implementation heavily relies on the author’s, I’m going to hope it’s the former.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t>After reviewing the instructions, I had to make a slight modification to my MCST function so that it returns a pointer to the MCST created, which can then be used to print it to a file.</w:t></w:r></w:p>
'@

$d = $word.ActiveDocument

function Get-ParaByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

function Get-ParaIndexByText($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# Apply edits from the BOTTOM of the document to the TOP so that
# paragraph-index shifts caused by inserted paragraphs never affect
# an edit that hasn't run yet.
# ------------------------------------------------------------------

# 1) "I've built the read, write, and print functions..." (+ 2 new paragraphs after)
$p44 = Get-ParaByText $d "I’ve built the read, write, and print functions"
$r44 = $p44.Range
$r44.InsertXML($block44)

# 2) "I'm gonna use a smart pointer to store the graph in main.cpp" (+ 3 new paragraphs after)
$p42 = Get-ParaByText $d "gonna use a smart pointer"
$r42 = $p42.Range
$r42.InsertXML($block42)

# 3) Insert new paragraph "Write function for mcst" immediately before "Write main"
#    (the "Write main" paragraph that follows "Write function to print graph to console")
$idxWriteMain = Get-ParaIndexByText $d "Write function to print graph to console"
$idxWriteMain = $idxWriteMain + 1
$pWriteMain = $d.Paragraphs($idxWriteMain)
$rWriteMain = $pWriteMain.Range
# Preserve "Write main" unchanged, right after our new paragraph (same pPr as before).
$writeMainBlock = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t>Write main</w:t></w:r></w:p>'
$combined = $blockWritefn + $writeMainBlock
$rWriteMain.InsertXML($combined)

# 4) "Write main that tests reading graph from file..." paragraph, plus removing one
#    of the two empty paragraphs that follow it.
$idxMcstMain = Get-ParaIndexByText $d "Write main that tests reading graph from file"
$pA = $d.Paragraphs($idxMcstMain)
$pB = $d.Paragraphs($idxMcstMain + 2)   # second (last) of the two empty paragraphs
$rangeAB = $d.Range($pA.Range.Start, $pB.Range.End)
$rangeAB.InsertXML($block29)

# 5) "Now that the GraphADT class has been built..." paragraph
$p16 = Get-ParaByText $d "has been built, I" 
$r16 = $p16.Range
$r16.InsertXML($block16)

# 6) "I'm going to start by building the GraphADT class..." paragraph
$p15 = Get-ParaByText $d "I’m going to start by building the GraphADT"
$r15 = $p15.Range
$r15.InsertXML($block15)

# 7) "I used a smart pointer to manage the graph object in main.cpp." (+ new paragraph after)
$p3 = Get-ParaByText $d "I used a smart pointer to manage the graph object"
$r3 = $p3.Range
$r3.InsertXML($block3)

Write-Output "All edits applied."
